$wb = $excel.ActiveWorkbook

# Sheet "展览" (worksheet 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 8174
$ws1.Range("F5").Value = 5955
$ws1.Range("F6").Value = 505
$ws1.Range("F7").Value = 96
$ws1.Range("F11").Value = 714
$ws1.Range("F12").Value = 72

# Sheet "全部类型" (worksheet 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 8174
$ws4.Range("F5").Value = 5955
$ws4.Range("F6").Value = 505
$ws4.Range("F7").Value = 96
$ws4.Range("F15").Value = 714
$ws4.Range("F16").Value = 72
